# This script reproduces, via Word COM-interop calls, the following
# textual corrections made to the "Seguimiento" tables in the document:
#   1) Two cells containing "Ninguna" become "Ningun" + "o" (two runs).
#   2) A cell containing the misspelling "Niguno" becomes "Ni" + "n" + "guno"
#      (three runs), and the paragraph mark picks up the run formatting
#      (i.e. the paragraph's <w:pPr> gains an <w:rPr> matching the runs).
#   3) A cell containing "Plazo al dia" becomes "Ninguno".

$d = $word.ActiveDocument

function Find-NextRange($searchText) {
    # Search from the top of the story each time: once a match is edited it
    # no longer matches the search text, so the next call naturally lands on
    # the following occurrence.
    $probe = $d.Content
    $found = $probe.Find.Execute($searchText, $true, $false, $false, $false, `
                                  $false, $true, 1, $false, "", 0)
    if (-not $found) {
        return $null
    }
    return $d.Range($probe.Start, $probe.End)
}

function Toggle-Bold($rng) {
    # Flipping Bold on and back off does not change the visible formatting
    # (both states resolve to "not bold"), but it forces the host to mint a
    # distinct run for $rng instead of silently re-merging it with an
    # adjacent run that happens to already share the same formatting.
    $rng.Bold = 1
    $rng.Bold = 0
}

# --- 1) & 2): the two "Ninguna" -> "Ningun" + "o" cells -------------------
for ($i = 0; $i -lt 2; $i++) {
    $r = Find-NextRange("Ninguna")
    if ($null -eq $r) { break }

    $start = $r.Start
    $end = $r.End

    # Replace the whole run's text with "Ningun" ...
    $r.Text = "Ningun"
    Toggle-Bold($r)

    # ... then append the trailing "o" as its own run, right after it.
    $r.Collapse(0)
    $r.InsertAfter("o")
    $tailStart = $start + 6
    $tail = $d.Range($tailStart, $tailStart + 1)
    Toggle-Bold($tail)
}

# --- 3): "Niguno" -> "Ni" + "n" + "guno" -----------------------------------
$r = Find-NextRange("Niguno")
if ($null -ne $r) {
    $start = $r.Start

    # Replace the whole run's text with "Ni" ...
    $r.Text = "Ni"
    # ... then append "n" and "guno" as their own trailing runs.
    $r.Collapse(0)
    $r.InsertAfter("n")
    $r.Collapse(0)
    $r.InsertAfter("guno")

    # Now split into three distinct runs with identical formatting, and
    # make sure every one of them (plus the paragraph mark, implicitly
    # carried along because the full original run's text got touched)
    # ends up with the same explicit formatting as before.
    $seg1 = $d.Range($start, $start + 2)      # "Ni"
    Toggle-Bold($seg1)
    $seg2 = $d.Range($start + 2, $start + 3)  # "n"
    Toggle-Bold($seg2)
    $seg3 = $d.Range($start + 3, $start + 7)  # "guno"
    Toggle-Bold($seg3)
}

# --- 4): "Plazo al dia" -> "Ninguno" ---------------------------------------
$r = Find-NextRange("Plazo al dia")
if ($null -ne $r) {
    $r.Text = "Ninguno"
}
